$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 351
$ws.Range("D6").Value = 142
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 3

$ws.Range("E20").Select()
